# appConfigFileMaryKay.xlsx - "New Code with API"
#
# 1. Remove the obsolete "Login" sheet (the suite now logs in as part of
#    the guest checkout flow instead of a dedicated Login page).
# 2. Add new locators for the search flow (SearchPage) and the guest
#    checkout heading / button (Checkout), add a PDP "go to checkout"
#    locator, fix a couple of typos, and point the Checkout ZipCode entry
#    at the already-existing "ZipCode" key.
# 3. Leave the window looking at the Checkout sheet / A9 like the author did.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Drop the Login sheet -------------------------------------------------
$wb.Worksheets.Item("Login").Delete()

$searchPage  = $wb.Worksheets.Item("SearchPage")
$pdp         = $wb.Worksheets.Item("PDP")
$checkout    = $wb.Worksheets.Item("Checkout")
$shipperForm = $wb.Worksheets.Item("ShipperForm")

# --- 2a. SearchPage -----------------------------------------------------------
# Insert two new rows (search / searchedProduct) after searchBar, pushing the
# old searchResult row down.
$searchPage.Range("A3").Value = "search"
$searchPage.Range("B3").Value = "(//*[@id='header-search-submit-button'])[1]"

$searchPage.Range("A4").Value = "searchResult"
$searchPage.Range("B4").Value = "(//*[@class='product-name'])[1]"

$searchPage.Range("A5").Value = "searchedProduct"
$searchPage.Range("B5").Value = "(//*[@class='product-name'])[1]"

$searchPage.Columns.Item(1).ColumnWidth = 15.166666666666666

# --- 2b. PDP --------------------------------------------------------------
# New row pointing at the "Checkout" link from the product page.
$pdp.Range("A7").Value = "check-Out"
$pdp.Range("B7").Value = "(//*[contains(text(),'Checkout')])[2]"

# --- 2c. Checkout -----------------------------------------------------------
# Fix casing on the zip-code key and the guest-checkout typo, then add the
# new guest-checkout heading + button locators.
$checkout.Range("A5").Value = "ZipCode"

$checkout.Range("A6").Value = "guestChek-Out"

$checkout.Range("A8").Value = "guestCheckOutHeading"
$checkout.Range("B8").Value = "(//h3)[2]"

$checkout.Range("A9").Value = "checkOutbtn"
$checkout.Range("B9").Value = "(//*[contains(text(),'Checkout')])[3]"

# --- 2d. Leave the cursor on each sheet where the author left it, then
#         activate Checkout / A9 last since that's the tab Excel shows on
#         reopen.
$shipperForm.Range("A10").Select()
$searchPage.Range("B12").Select()
$pdp.Range("B6").Select()

$checkout.Activate()
$checkout.Range("A9").Select()

Write-Output "appConfigFileMaryKay.xlsx updated"
